# Suzie's IRA workbook update: record a new AT&T (T) dividend/reinvestment
# entry on the "Initial Buys" sheet, and update the corresponding share
# count / reinvestment figures on the "2016" sheet. Everything else
# (2017 totals, 2016/2017 sum rows, etc.) is formula-driven and will
# recalculate automatically.

$wb = $excel.ActiveWorkbook

# ---- Initial Buys ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Initial Buys")

# New dividend entry for ticker T (AT&T): date + amount in row 3,
# matching the style already used by the row above (row 2).
$ws1.Range("V2").Copy()
$ws1.Range("V3").PasteSpecial(-4122)
$ws1.Range("W2").Copy()
$ws1.Range("W3").PasteSpecial(-4122)

$ws1.Range("V3").Value = 42667
$ws1.Range("W3").Value = 654.66

[void]$ws1.Range("M18").Select()

# ---- 2016 -------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("2016")

# PM (row 7): dividend was reinvested, bumping share count and recording
# the reinvestment amount in the Oct-2016 column (Q).
$ws3.Range("D7").Value = 5.054
$ws3.Range("Q7").Value = 5.2

# T (row 9): share count increased after the new Initial Buys entry above
# was reinvested.
$ws3.Range("D9").Value = 45

[void]$ws3.Range("N16").Select()
